# CW2 development, 1st section and user account section complete
#
# 1) Bump the cached "datetimeFigureOut" text (Date Placeholder) from
#    21/07/2025 -> 22/07/2025 on the slide master and every slide layout.
# 2) Slide 15 ("groups"): move "Picture 4" down and delete the empty
#    "Content Placeholder 8" shape.
# 3) Slide 16: nudge "Picture 6" to the left.

$p = $ppt.ActivePresentation

$oldDate = "21/07/2025"
$newDate = "22/07/2025"

function Update-DateShape($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# Slide master date placeholder.
$master = $p.SlideMaster
Update-DateShape $master.Shapes

# Every slide layout (custom layout) date placeholder.
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DateShape $layout.Shapes
}

# --- Slide 15: "groups" ---
$slide15 = $p.Slides.Item(15)
for ($i = 1; $i -le $slide15.Shapes.Count; $i++) {
    $shp = $slide15.Shapes.Item($i)
    if ($shp.Name -eq "Picture 4") {
        $shp.Top = 1411589 / 12700
    }
}

for ($i = $slide15.Shapes.Count; $i -ge 1; $i--) {
    $shp = $slide15.Shapes.Item($i)
    if ($shp.Name -eq "Content Placeholder 8") {
        $shp.Delete()
    }
}

# --- Slide 16 ---
$slide16 = $p.Slides.Item(16)
for ($i = 1; $i -le $slide16.Shapes.Count; $i++) {
    $shp = $slide16.Shapes.Item($i)
    if ($shp.Name -eq "Picture 6") {
        $shp.Left = 162560 / 12700
    }
}
